# Merge the split runs in the Title / Author / Abstract paragraphs into a
# single run each (collapsing "Factsheet:" + " " + "Lognormal" + ... into
# one run of text), matching the canonical re-saved OOXML.
#
# Word's object model has no direct "Runs" collection, so each paragraph is
# rebuilt the way a VBA macro would: keep the paragraph's first run (so its
# xml:space="preserve" formatting sticks), wipe everything after it, then
# re-insert the remaining words with InsertAfter so the whole paragraph
# collapses back down to a single <w:r>.

$d = $word.ActiveDocument

function Merge-ParagraphRuns {
    param(
        [int]$ParaIndex,
        [string]$FullText
    )

    $p = $d.Paragraphs.Item($ParaIndex)
    $paraStart = $p.Range.Start
    $paraEnd = $p.Range.End - 1   # exclude the paragraph mark

    # Keep the very first run/word intact (preserves its xml:space attr),
    # then clear everything else in the paragraph and retype it onto the
    # end of that first run so it all collapses into one run.
    $firstSpace = $FullText.IndexOf(" ")
    if ($firstSpace -lt 0) {
        $firstWord = $FullText
    } else {
        $firstWord = $FullText.Substring(0, $firstSpace)
    }
    $remainder = $FullText.Substring($firstWord.Length)

    $splitPoint = $paraStart + $firstWord.Length

    if ($paraEnd -gt $splitPoint) {
        $tail = $d.Range($splitPoint, $paraEnd)
        $tail.Text = ""
    }

    if ($remainder.Length -gt 0) {
        $insertPoint = $d.Range($splitPoint, $splitPoint)
        $insertPoint.InsertAfter($remainder)
    }
}

Merge-ParagraphRuns 1 "Factsheet: Lognormal distribution"
Merge-ParagraphRuns 2 "Michelle Arnetta and Tom Coleman"
Merge-ParagraphRuns 4 "A factsheet for the lognormal distribution."
